$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Swap-Row($ws, $row1, $row2) {
    # Swap columns B and F..AC (everything except A, C, D, E which are identical
    # between the two matches already) between the two given rows.
    $cols = @(2,6,7,8,9,10,11,12,13,14,15,16,17,18,19,20,21,22,23,24,25,26,27,28,29)
    foreach ($col in $cols) {
        $c1 = $ws.Cells.Item($row1, $col)
        $c2 = $ws.Cells.Item($row2, $col)
        $v1 = $c1.Value()
        $v2 = $c2.Value()
        $c1.Value = $v2
        $c2.Value = $v1
    }
}

# Rows 114/115 (ids 112/113) had their B..AC (minus C/D/E) data mixed up between
# the two matches; fix by swapping the two rows' content back.
Swap-Row $ws 114 115

# Same issue for rows 119/120 (ids 117/118).
Swap-Row $ws 119 120

# Append the new match result as row 177 (id 175). Copy row 176's formatting
# first so the new row reuses the existing cell styles (bold/boxed id column,
# date-formatted date column) instead of creating new style entries.
$ws.Range("A176:AC176").Copy($ws.Range("A177:AC177"))

# Clear the copied values for columns that should stay empty on this new row
# (no result yet: FTHG/FTAG/FTR, and no closing-line PL_AhOver/PL_AhUnder yet).
$ws.Range("H177:J177").ClearContents()
$ws.Range("AB177:AC177").ClearContents()

$ws.Cells.Item(177, 1).Value = 175
$ws.Cells.Item(177, 2).Value = 8051187
$ws.Cells.Item(177, 3).Value = "Uruguay Primera División"
$ws.Cells.Item(177, 4).Value = "Uruguay Apertura"
$ws.Cells.Item(177, 5).Value = 45396.625
$ws.Cells.Item(177, 6).Value = "Defensor Sporting"
$ws.Cells.Item(177, 7).Value = "CA River Plate"
$ws.Cells.Item(177, 11).Value = 1.727
$ws.Cells.Item(177, 12).Value = 3.5
$ws.Cells.Item(177, 13).Value = 5
$ws.Cells.Item(177, 14).Value = 1.571
$ws.Cells.Item(177, 15).Value = 3.6
$ws.Cells.Item(177, 16).Value = 6.5
$ws.Cells.Item(177, 17).Value = -1
$ws.Cells.Item(177, 18).Value = 2.025
$ws.Cells.Item(177, 19).Value = 1.825
$ws.Cells.Item(177, 20).Value = 2.25
$ws.Cells.Item(177, 21).Value = 1.8
$ws.Cells.Item(177, 22).Value = 2.05
$ws.Cells.Item(177, 23).Value = 0
$ws.Cells.Item(177, 24).Value = 0
$ws.Cells.Item(177, 25).Value = 0
$ws.Cells.Item(177, 26).Value = 0
$ws.Cells.Item(177, 27).Value = 0
